$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new "Từ khóa:" label in D5
$ws.Range("D5").Value = "Từ khóa:"

# Delete rows 11:12 (empty spacer rows), shifting rows below up by two
$ws.Rows("11:12").Delete()

# Re-anchor the frozen pane to the new header row (was row 14/15, now 12/13)
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A13").Select()
$excel.ActiveWindow.FreezePanes = $true

# Update selection to match the target state
$ws.Range("D15").Select()
